$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 26 and row 27 swap places (Kaspa <-> Cosmos) along with their updated
# price/volume figures. Column A (rank index) stays fixed per row; only the
# coin name, link, price and volume columns change content.
$ws.Range('B26').Value = 'Cosmos'
$ws.Range('C26').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$c = $ws.Range('D26')
$c.NumberFormat = '@'
$c.Value = '9.16'
$c.Style = 'Normal'
$ws.Range('E26').Value = '  -1.18%  '

$ws.Range('B27').Value = 'Kaspa'
$ws.Range('C27').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$c = $ws.Range('D27')
$c.NumberFormat = '@'
$c.Value = '0.138'
$c.Style = 'Normal'
$ws.Range('E27').Value = '  +6.51%  '

# Updated price / 1h-volume figures for the remaining coin rows. Cells that
# hold a plain-numeric-looking price (e.g. "244.73") are forced to text via
# a temporary "@" number format so Excel doesn't silently convert them to
# numeric cells (the source data stores these as text), then the format is
# reset back to Normal so no stray cell formatting is introduced.
$ws.Range('D2').Value = '36.598.85'
$ws.Range('E2').Value = '  +0.74%  '
$ws.Range('D3').Value = '1.960.94'
$ws.Range('E3').Value = '  +1.53%  '
$ws.Range('E4').Value = '  -0.09%  '
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '244.73'
$c.Style = 'Normal'
$ws.Range('E5').Value = '  +1.55%  '
$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '0.616'
$c.Style = 'Normal'
$ws.Range('E6').Value = '  +1.55%  '
$c = $ws.Range('D7')
$c.NumberFormat = '@'
$c.Value = '58.78'
$c.Style = 'Normal'
$ws.Range('E7').Value = '  +3.54%  '
$ws.Range('E8').Value = '  -0.04%  '
$ws.Range('E9').Value = '  +4.70%  '
$c = $ws.Range('D10')
$c.NumberFormat = '@'
$c.Value = '0.0793'
$c.Style = 'Normal'
$ws.Range('E10').Value = '  -4.90%  '
$ws.Range('E11').Value = '  -1.04%  '
$c = $ws.Range('D12')
$c.NumberFormat = '@'
$c.Value = '14.22'
$c.Style = 'Normal'
$ws.Range('E12').Value = '  +6.35%  '
$c = $ws.Range('D13')
$c.NumberFormat = '@'
$c.Value = '0.838'
$c.Style = 'Normal'
$ws.Range('E13').Value = '  +4.66%  '
$ws.Range('D14').Value = '2.246.53'
$ws.Range('E14').Value = '  +1.37%  '
$c = $ws.Range('D15')
$c.NumberFormat = '@'
$c.Value = '21.19'
$c.Style = 'Normal'
$ws.Range('E15').Value = '  +1.43%  '
$c = $ws.Range('D16')
$c.NumberFormat = '@'
$c.Value = '5.30'
$c.Style = 'Normal'
$ws.Range('E16').Value = '  +3.32%  '
$ws.Range('D17').Value = '1.955.70'
$ws.Range('E17').Value = '  +0.12%  '
$ws.Range('D18').Value = '36.535.35'
$ws.Range('E18').Value = '  +0.76%  '
$c = $ws.Range('D19')
$c.NumberFormat = '@'
$c.Value = '69.88'
$c.Style = 'Normal'
$ws.Range('E19').Value = '  +1.56%  '
$ws.Range('D20').Value = '0.0₃0850'
$ws.Range('E20').Value = '  -0.97%  '
$c = $ws.Range('D21')
$c.NumberFormat = '@'
$c.Value = '229.87'
$c.Style = 'Normal'
$ws.Range('E21').Value = '  +1.87%  '
$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value = '5.06'
$c.Style = 'Normal'
$ws.Range('E22').Value = '  +2.57%  '
$c = $ws.Range('D23')
$c.NumberFormat = '@'
$c.Value = '0.999'
$c.Style = 'Normal'
$ws.Range('E23').Value = '  -0.07%  '
$ws.Range('E24').Value = '  +5.87%  '
$ws.Range('E25').Value = '  +3.67%  '
$c = $ws.Range('D28')
$c.NumberFormat = '@'
$c.Value = '160.84'
$c.Style = 'Normal'
$ws.Range('E28').Value = '  +0.31%  '
$c = $ws.Range('D29')
$c.NumberFormat = '@'
$c.Value = '19.45'
$c.Style = 'Normal'
$ws.Range('E29').Value = '  +1.87%  '
$c = $ws.Range('D30')
$c.NumberFormat = '@'
$c.Value = '1.21'
$c.Style = 'Normal'
$ws.Range('E30').Value = '  +9.23%  '
$ws.Range('E31').Value = '  +2.44%  '
$c = $ws.Range('D32')
$c.NumberFormat = '@'
$c.Value = '4.73'
$c.Style = 'Normal'
$ws.Range('E32').Value = '  +4.44%  '
$c = $ws.Range('D33')
$c.NumberFormat = '@'
$c.Value = '0.0613'
$c.Style = 'Normal'
$ws.Range('E33').Value = '  -1.93%  '
$c = $ws.Range('D34')
$c.NumberFormat = '@'
$c.Value = '4.40'
$c.Style = 'Normal'
$ws.Range('E34').Value = '  +6.37%  '
$ws.Range('E35').Value = '  +18.70%  '
$c = $ws.Range('D36')
$c.NumberFormat = '@'
$c.Value = '2.30'
$c.Style = 'Normal'
$ws.Range('E36').Value = '  +8.66%  '
$ws.Range('E37').Value = '  -0.22%  '
$c = $ws.Range('D38')
$c.NumberFormat = '@'
$c.Value = '1.77'
$c.Style = 'Normal'
$ws.Range('E38').Value = '  -0.99%  '
$c = $ws.Range('D39')
$c.NumberFormat = '@'
$c.Value = '5.46'
$c.Style = 'Normal'
$ws.Range('E39').Value = '  -10.18%  '
$c = $ws.Range('D40')
$c.NumberFormat = '@'
$c.Value = '0.0984'
$c.Style = 'Normal'
$ws.Range('E40').Value = '  +1.29%  '
$ws.Range('E41').Value = '  +1.43%  '
$ws.Range('E42').Value = '  +2.21%  '
$ws.Range('E43').Value = '  +1.51%  '
$ws.Range('D44').Value = '1.375.85'
$ws.Range('E44').Value = '  +3.30%  '
$c = $ws.Range('D45')
$c.NumberFormat = '@'
$c.Value = '15.81'
$c.Style = 'Normal'
$ws.Range('E45').Value = '  +2.34%  '
$c = $ws.Range('D46')
$c.NumberFormat = '@'
$c.Value = '88.12'
$c.Style = 'Normal'
$ws.Range('E46').Value = '  +2.38%  '
$c = $ws.Range('D47')
$c.NumberFormat = '@'
$c.Value = '1.03'
$c.Style = 'Normal'
$ws.Range('E47').Value = '  +1.49%  '
$c = $ws.Range('D48')
$c.NumberFormat = '@'
$c.Value = '7.14'
$c.Style = 'Normal'
$ws.Range('E48').Value = '  +1.35%  '
$ws.Range('E49').Value = '  +0.84%  '
$ws.Range('D50').Value = '2.137.02'
$ws.Range('E50').Value = '  +1.37%  '
$c = $ws.Range('D51')
$c.NumberFormat = '@'
$c.Value = '44.02'
$c.Style = 'Normal'
$ws.Range('E51').Value = '  +0.55%  '
